$wb = $excel.ActiveWorkbook

# This script applies a scheduled market-data refresh to the Tiamat_Profits
# workbook: updated currentAveragePrice / LevePrice / LeveProfit figures
# (columns H-N) for the affected leve rows, across all 8 job sheets.

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 514257.16
$ws.Range("J17").Value = 514257.16
$ws.Range("L17").Value = 1542771.48
$ws.Range("N17").Value = -1543107.48
$ws.Range("H98").Value = 672.96875
$ws.Range("I98").Value = 538.4815
$ws.Range("J98").Value = 1399.2
$ws.Range("K98").Value = 538.4815
$ws.Range("L98").Value = 1399.2
$ws.Range("M98").Value = 959.5185
$ws.Range("N98").Value = -4395.2
$ws.Range("H122").Value = 672.96875
$ws.Range("I122").Value = 538.4815
$ws.Range("J122").Value = 1399.2
$ws.Range("K122").Value = 1615.4445
$ws.Range("L122").Value = 4197.6
$ws.Range("M122").Value = 834.5554999999999
$ws.Range("N122").Value = -9097.6
$ws.Range("H127").Value = 806
$ws.Range("I127").Value = 333.1
$ws.Range("J127").Value = 1020.9545
$ws.Range("K127").Value = 999.3000000000001
$ws.Range("L127").Value = 3062.8635
$ws.Range("M127").Value = 3960.7
$ws.Range("N127").Value = -12982.8635
$ws.Range("H129").Value = 2114.5417
$ws.Range("I129").Value = 353.6111
$ws.Range("J129").Value = 7397.3335
$ws.Range("K129").Value = 1060.8333
$ws.Range("L129").Value = 22192.0005
$ws.Range("M129").Value = 3939.1667
$ws.Range("N129").Value = -32192.0005
$ws.Range("H135").Value = 25001110
$ws.Range("I135").Value = 318.07693
$ws.Range("K135").Value = 2862.69237
$ws.Range("M135").Value = -327.6923700000002

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2266.6667
$ws.Range("J63").Value = 1800
$ws.Range("L63").Value = 1800
$ws.Range("N63").Value = -3172
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H66").Value = 2266.6667
$ws.Range("J66").Value = 1800
$ws.Range("L66").Value = 9000
$ws.Range("N66").Value = -15864
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H86").Value = 111127780
$ws.Range("J86").Value = 142875710
$ws.Range("L86").Value = 142875710
$ws.Range("N86").Value = -142878082
$ws.Range("H89").Value = 111127780
$ws.Range("J89").Value = 142875710
$ws.Range("L89").Value = 428627130
$ws.Range("N89").Value = -428638986

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 390717.16
$ws.Range("I86").Value = 1881.091
$ws.Range("J86").Value = 1001745.3
$ws.Range("K86").Value = 1881.091
$ws.Range("L86").Value = 1001745.3
$ws.Range("M86").Value = -758.0909999999999
$ws.Range("N86").Value = -1003991.3
$ws.Range("H89").Value = 390717.16
$ws.Range("I89").Value = 1881.091
$ws.Range("J89").Value = 1001745.3
$ws.Range("K89").Value = 9405.455
$ws.Range("L89").Value = 5008726.5
$ws.Range("M89").Value = -3789.455
$ws.Range("N89").Value = -5019958.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38023.49
$ws.Range("I31").Value = 67870.05499999999
$ws.Range("J31").Value = 9747.789000000001
$ws.Range("K31").Value = 67870.05499999999
$ws.Range("L31").Value = 9747.789000000001
$ws.Range("M31").Value = -67575.05499999999
$ws.Range("N31").Value = -10337.789
$ws.Range("H34").Value = 38023.49
$ws.Range("I34").Value = 67870.05499999999
$ws.Range("J34").Value = 9747.789000000001
$ws.Range("K34").Value = 67870.05499999999
$ws.Range("L34").Value = 9747.789000000001
$ws.Range("M34").Value = -67668.05499999999
$ws.Range("N34").Value = -10151.789
$ws.Range("H58").Value = 3502.8333
$ws.Range("I58").Value = 1242.3334
$ws.Range("J58").Value = 5763.3335
$ws.Range("K58").Value = 1242.3334
$ws.Range("L58").Value = 5763.3335
$ws.Range("M58").Value = -1039.3334
$ws.Range("N58").Value = -6169.3335
$ws.Range("H62").Value = 25002366
$ws.Range("I62").Value = 45457140
$ws.Range("J62").Value = 2088.889
$ws.Range("K62").Value = 45457140
$ws.Range("L62").Value = 2088.889
$ws.Range("M62").Value = -45456516
$ws.Range("N62").Value = -3336.889
$ws.Range("H65").Value = 25002366
$ws.Range("I65").Value = 45457140
$ws.Range("J65").Value = 2088.889
$ws.Range("K65").Value = 227285700
$ws.Range("L65").Value = 10444.445
$ws.Range("M65").Value = -227282580
$ws.Range("N65").Value = -16684.445
$ws.Range("H132").Value = 2235.7673
$ws.Range("I132").Value = 1835.2632
$ws.Range("J132").Value = 5279.6
$ws.Range("K132").Value = 5505.7896
$ws.Range("L132").Value = 15838.8
$ws.Range("M132").Value = -2975.7896
$ws.Range("N132").Value = -20898.8
$ws.Range("H134").Value = 33335508
$ws.Range("I134").Value = 1662.4
$ws.Range("J134").Value = 100003200
$ws.Range("K134").Value = 4987.200000000001
$ws.Range("L134").Value = 300009600
$ws.Range("M134").Value = -2452.200000000001
$ws.Range("N134").Value = -300014670
$ws.Range("H136").Value = 3502.8333
$ws.Range("I136").Value = 1242.3334
$ws.Range("J136").Value = 5763.3335
$ws.Range("K136").Value = 3727.0002
$ws.Range("L136").Value = 17290.0005
$ws.Range("M136").Value = -1177.0002
$ws.Range("N136").Value = -22390.0005

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 400000320
$ws.Range("I26").Value = 666666700
$ws.Range("J26").Value = 800
$ws.Range("K26").Value = 2000000100
$ws.Range("L26").Value = 2400
$ws.Range("M26").Value = -1999999812
$ws.Range("N26").Value = -2976
$ws.Range("H64").Value = 2168666.8
$ws.Range("J64").Value = 2168666.8
$ws.Range("L64").Value = 6506000.399999999
$ws.Range("N64").Value = -6506540.399999999
$ws.Range("H67").Value = 2168666.8
$ws.Range("J67").Value = 2168666.8
$ws.Range("L67").Value = 6506000.399999999
$ws.Range("N67").Value = -6507872.399999999
$ws.Range("H86").Value = 302
$ws.Range("I86").Value = 302
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 906
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 280
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 302
$ws.Range("I89").Value = 302
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 2718
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 3210
$ws.Range("N89").ClearContents()
$ws.Range("H134").Value = 4270.488
$ws.Range("I134").Value = 1296.1875
$ws.Range("J134").Value = 6174.04
$ws.Range("K134").Value = 3888.5625
$ws.Range("L134").Value = 18522.12
$ws.Range("M134").Value = 1181.4375
$ws.Range("N134").Value = -28662.12

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2162.3333
$ws.Range("I68").Value = 1708.4166
$ws.Range("J68").Value = 3070.1667
$ws.Range("K68").Value = 1708.4166
$ws.Range("L68").Value = 3070.1667
$ws.Range("M68").Value = -959.4166
$ws.Range("N68").Value = -4568.1667
$ws.Range("H71").Value = 2162.3333
$ws.Range("I71").Value = 1708.4166
$ws.Range("J71").Value = 3070.1667
$ws.Range("K71").Value = 8542.083000000001
$ws.Range("L71").Value = 15350.8335
$ws.Range("M71").Value = -4798.083000000001
$ws.Range("N71").Value = -22838.8335

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H87").Value = 15000
$ws.Range("J87").Value = 15000
$ws.Range("L87").Value = 15000
$ws.Range("N87").Value = -17496
$ws.Range("H90").Value = 15000
$ws.Range("J90").Value = 15000
$ws.Range("L90").Value = 45000
$ws.Range("N90").Value = -57480
$ws.Range("H132").Value = 3161.54
$ws.Range("I132").Value = 873.6923
$ws.Range("J132").Value = 11273
$ws.Range("K132").Value = 2621.0769
$ws.Range("L132").Value = 33819
$ws.Range("M132").Value = -91.07690000000002
$ws.Range("N132").Value = -38879
